# Updated cryptos list on Wed Jun  5 13:39:36 UTC 2024 with GitHub Actions
# Refreshes the scraped Price / Volume(1h) figures for every coin row, and
# replaces the Bittensor row with Cosmos (name, link, price, volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as plain text in this sheet (e.g.
# thousand-separator dots like "71.095.52", or trailing zeros like "700.10").
# Any new value that looks like an ordinary number needs its cell
# pre-formatted as Text so Excel keeps it as a literal string instead of
# silently converting/rounding it to a numeric value.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "70.990.76"
$ws.Range("E2").Value = "  +2.39%  "

# Row 3
$ws.Range("D3").Value = "3.807.29"

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "700.10"
$ws.Range("E5").Value = "  +8.46%  "

# Row 6
$ws.Range("D6").Value = "173.42"
$ws.Range("E6").Value = "  +4.48%  "

# Row 7
$ws.Range("D7").Value = "3.806.94"
$ws.Range("E7").Value = "  +0.94%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("E10").Value = "  +2.44%  "

# Row 11
$ws.Range("D11").Value = "7.24"
$ws.Range("E11").Value = "  +4.92%  "

# Row 12
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("E13").Value = "  +8.17%  "

# Row 14
$ws.Range("D14").Value = "36.22"
$ws.Range("E14").Value = "  +3.88%  "

# Row 15
$ws.Range("D15").Value = "4.450.88"
$ws.Range("E15").Value = "  +0.97%  "

# Row 16
$ws.Range("D16").Value = "3.807.76"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").Value = "71.010.72"
$ws.Range("E17").Value = "  +2.41%  "

# Row 18
$ws.Range("D18").Value = "17.78"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").Value = "7.19"
$ws.Range("E19").Value = "  +2.46%  "

# Row 20
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("E21").Value = "  +16.81%  "

# Row 22
$ws.Range("D22").Value = "479.07"
$ws.Range("E22").Value = "  +2.44%  "

# Row 23
$ws.Range("D23").Value = "0.711"
$ws.Range("E23").Value = "  +0.44%  "

# Row 24
$ws.Range("D24").Value = "83.83"
$ws.Range("E24").Value = "  +2.44%  "

# Row 25
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "2.16"
$ws.Range("E27").Value = "  +2.56%  "

# Row 28
$ws.Range("D28").Value = "10.42"
$ws.Range("E28").Value = "  +0.32%  "

# Row 29
$ws.Range("D29").Value = "3.961.08"
$ws.Range("E29").Value = "  +1.00%  "

# Row 30
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("E31").Value = "  +15.67%  "

# Row 32
$ws.Range("D32").Value = "2.29"
$ws.Range("E32").Value = "  +1.12%  "

# Row 33
$ws.Range("D33").Value = "7.51"
$ws.Range("E33").Value = "  +4.83%  "

# Row 34
$ws.Range("E34").Value = "  +8.52%  "

# Row 35
$ws.Range("D35").Value = "29.49"
$ws.Range("E35").Value = "  +2.83%  "

# Row 36
$ws.Range("D36").Value = "9.24"
$ws.Range("E36").Value = "  +4.34%  "

# Row 37
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("E38").Value = "  +2.54%  "

# Row 39
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +5.35%  "

# Row 40
$ws.Range("D40").Value = "6.00"
$ws.Range("E40").Value = "  +2.74%  "

# Row 41
$ws.Range("E41").Value = "  +13.52%  "

# Row 42
$ws.Range("D42").Value = "0.981"
$ws.Range("E42").Value = "  +2.41%  "

# Row 43
$ws.Range("D43").Value = "0.000328"
$ws.Range("E43").Value = "  +21.41%  "

# Row 44
$ws.Range("E44").Value = "  -0.10%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").Value = "163.14"
$ws.Range("E46").Value = "  +4.14%  "

# Row 47
$ws.Range("D47").Value = "48.97"
$ws.Range("E47").Value = "  +3.13%  "

# Row 48
$ws.Range("D48").Value = "44.34"
$ws.Range("E48").Value = "  -2.14%  "

# Row 49
$ws.Range("D49").Value = "0.300"
$ws.Range("E49").Value = "  +0.97%  "

# Row 50
$ws.Range("D50").Value = "1.38"
$ws.Range("E50").Value = "  -1.52%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "8.56"
$ws.Range("E51").Value = "  +2.20%  "

